# Regenerate the quadratic/linear problem data (as in the commit
# "volver a generar problemas cuadraticos y lineales").
#
# All the touched cells hold plain text in the workbook (General-formatted,
# shared-string) values, even the ones that look like numbers. Writing a
# bare numeric-looking string via .Value lets Excel coerce it to a real
# number (and reformat/round it), so any such value is written with a
# leading apostrophe to force a text literal and preserve the exact digits.
# Cells whose text is not purely numeric (e.g. "-12.85 + x + 2y") are
# written without the apostrophe since Excel already keeps those as text.

# NOTE: sheet-name lookups (Worksheets.Item("Vector_bf") /
# Worksheets.Item("Vector_BF")) are case-insensitive, and this workbook has
# two sheets whose names differ only by case ("Vector_bf" vs "Vector_BF").
# To avoid ambiguity, address sheets by their (1-based) tab position instead:
#   1 Funciones_Objetivo        5 Vector_bf
#   2 Restricciones_del_lider   6 Vector_BF
#   3 Restricciones_del_follower 7 Vector_Alpha
#   4 Punto_modificado

$wb = $excel.ActiveWorkbook

$wsFollower = $wb.Worksheets.Item(3)   # Restricciones_del_follower
$wsFollower.Range("A2").Value = "-12.85 + x + 2y"
$wsFollower.Range("B2").Value = "'-1.1500000000000004"
$wsFollower.Range("D2").Value = "'0.07"
$wsFollower.Range("E2").Value = "'0.8999999999999999"
$wsFollower.Range("F2").Value = "'5.699999999999999"

$wsFollower.Range("A3").Value = "3.3499999999999996 + x - 2y"
$wsFollower.Range("B3").Value = "'-5.35"
$wsFollower.Range("D3").Value = "'0.21"
$wsFollower.Range("E3").Value = "'5.300000000000001"
$wsFollower.Range("F3").Value = "'6.4"

$wsFollower.Range("A4").Value = "-7.45 - 2x + y"
$wsFollower.Range("B4").Value = "'-6.45"
$wsFollower.Range("D4").Value = "'0.4"
$wsFollower.Range("E4").Value = "'9.7"
$wsFollower.Range("F4").Value = "'2.1"

$wsPunto = $wb.Worksheets.Item(4)   # Punto_modificado
$wsPunto.Range("A2").Value = "'4.75"
$wsPunto.Range("B2").Value = "'4.05"

$wsVecBf = $wb.Worksheets.Item(5)   # Vector_bf
$wsVecBf.Range("A2").Value = "'1.7800000000000002"

$wsVecBF = $wb.Worksheets.Item(6)   # Vector_BF
$wsVecBF.Range("A2").Value = "'9.699999999999998"
$wsVecBF.Range("A3").Value = "'-4.999999999999998"
